$wb = $excel.ActiveWorkbook

# New timestamp values (recorded as plain text, same format as existing Date column)
$sheet1Dates = @(
    "Mon Oct 09 22:24:57 EDT 2023",
    "Mon Oct 09 22:25:09 EDT 2023",
    "Mon Oct 09 22:25:23 EDT 2023",
    "Mon Oct 09 22:25:35 EDT 2023",
    "Mon Oct 09 22:25:47 EDT 2023",
    "Mon Oct 09 22:25:59 EDT 2023",
    "Mon Oct 09 22:26:11 EDT 2023",
    "Mon Oct 09 22:26:22 EDT 2023",
    "Mon Oct 09 22:26:34 EDT 2023",
    "Mon Oct 09 22:26:46 EDT 2023",
    "Mon Oct 09 22:26:58 EDT 2023",
    "Mon Oct 09 22:27:09 EDT 2023"
)

$sheet2Dates = @(
    "Mon Oct 09 22:27:22 EDT 2023",
    "Mon Oct 09 22:27:33 EDT 2023",
    "Mon Oct 09 22:27:44 EDT 2023",
    "Mon Oct 09 22:27:55 EDT 2023",
    "Mon Oct 09 22:28:06 EDT 2023",
    "Mon Oct 09 22:28:18 EDT 2023",
    "Mon Oct 09 22:28:29 EDT 2023",
    "Mon Oct 09 22:28:40 EDT 2023",
    "Mon Oct 09 22:28:51 EDT 2023",
    "Mon Oct 09 22:29:02 EDT 2023",
    "Mon Oct 09 22:29:13 EDT 2023",
    "Mon Oct 09 22:29:25 EDT 2023",
    "Mon Oct 09 22:29:35 EDT 2023",
    "Mon Oct 09 22:29:47 EDT 2023",
    "Mon Oct 09 22:29:58 EDT 2023",
    "Mon Oct 09 22:30:09 EDT 2023"
)

$ws1 = $wb.Worksheets.Item("FEINmismatch")
for ($i = 0; $i -lt $sheet1Dates.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $sheet1Dates[$i]
}

$ws2 = $wb.Worksheets.Item("FEINSSNmismatch")
for ($i = 0; $i -lt $sheet2Dates.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 2).Value = $sheet2Dates[$i]
}
